$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header rename (column names) ---
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# --- Title-case municipality/state names (de/del/el/la/las/los/y -> De/Del/El/La/Las/Los/Y) ---
$ws.Range('B6').Value = 'Pabellón De Arteaga'
$ws.Range('B7').Value = 'Rincón De Romos'
$ws.Range('B38').Value = 'San Cristóbal De Las Casas'
$ws.Range('B63').Value = 'Guadalupe Y Calvo'
$ws.Range('B65').Value = 'Hidalgo Del Parral'
$ws.Range('B81').Value = 'San Francisco Del Oro'
$ws.Range('B85').Value = 'Valle De Zaragoza'
$ws.Range('B97').Value = 'San Juan De Sabinas'
$ws.Range('A110').Value = 'Ciudad De México'
$ws.Range('B114').Value = 'Cuajimalpa De Morelos'
$ws.Range('B126').Value = 'Coneto De Comonfort'
$ws.Range('B139').Value = 'Nombre De Dios'
$ws.Range('B142').Value = 'Pánuco De Coronado'
$ws.Range('B147').Value = 'San Juan Del Río'
$ws.Range('A155').Value = 'Estado De México'
$ws.Range('B155').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B157').Value = 'Almoloya De Juárez'
$ws.Range('B161').Value = 'Atizapán De Zaragoza'
$ws.Range('B166').Value = 'Coacalco De Berriozábal'
$ws.Range('B169').Value = 'Ecatepec De Morelos'
$ws.Range('B172').Value = 'Ixtapan De La Sal'
$ws.Range('B178').Value = 'Naucalpan De Juárez'
$ws.Range('B183').Value = 'San Felipe Del Progreso'
$ws.Range('B184').Value = 'San Martín De Las Pirámides'
$ws.Range('B189').Value = 'Tlalnepantla De Baz'
$ws.Range('B195').Value = 'Valle De Bravo'
$ws.Range('B196').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B197').Value = 'Villa Del Carbón'
$ws.Range('B203').Value = 'San Miguel De Allende'
$ws.Range('B204').Value = 'Apaseo El Alto'
$ws.Range('B205').Value = 'Apaseo El Grande'
$ws.Range('B212').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B216').Value = 'Jaral Del Progreso'
$ws.Range('B223').Value = 'Purísima Del Rincón'
$ws.Range('B227').Value = 'San Diego De La Unión'
$ws.Range('B229').Value = 'San Francisco Del Rincón'
$ws.Range('B231').Value = 'San Luis De La Paz'
$ws.Range('B232').Value = 'Santa Cruz De Juventino Rosas'
$ws.Range('B233').Value = 'Silao De La Victoria'
$ws.Range('B238').Value = 'Valle De Santiago'
$ws.Range('B243').Value = 'Acapulco De Juárez'
$ws.Range('B244').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B247').Value = 'Atenango Del Río'
$ws.Range('B249').Value = 'Ayutla De Los Libres'
$ws.Range('B252').Value = 'Chilapa De Álvarez'
$ws.Range('B253').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B256').Value = 'Coyuca De Benítez'
$ws.Range('B257').Value = 'Coyuca De Catalán'
$ws.Range('B261').Value = 'Cutzamala De Pinzón'
$ws.Range('B266').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B267').Value = 'Iguala De La Independencia'
$ws.Range('B268').Value = 'Ixcateopan De Cuauhtémoc'
$ws.Range('B271').Value = 'La Unión De Isidoro Montes De Oca'
$ws.Range('B281').Value = 'Taxco De Alarcón'
$ws.Range('B283').Value = 'Técpan De Galeana'
$ws.Range('B285').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B287').Value = 'Tixtla De Guerrero'
$ws.Range('B289').Value = 'Tlapa De Comonfort'
$ws.Range('B297').Value = 'Atotonilco El Grande'
$ws.Range('B299').Value = 'Cuautepec De Hinojosa'
$ws.Range('B304').Value = 'Jacala De Ledezma'
$ws.Range('B306').Value = 'Mixquiahuala De Juárez'
$ws.Range('B308').Value = 'Omitlán De Juárez'
$ws.Range('B309').Value = 'Pachuca De Soto'
$ws.Range('B311').Value = 'Progreso De Obregón'
$ws.Range('B317').Value = 'Tepehuacán De Guerrero'
$ws.Range('B318').Value = 'Tezontepec De Aldama'
$ws.Range('B320').Value = 'Tula De Allende'
$ws.Range('B321').Value = 'Tulancingo De Bravo'
$ws.Range('B325').Value = 'Ahualulco De Mercado'
$ws.Range('B328').Value = 'Atemajac De Brizuela'
$ws.Range('B330').Value = 'Atotonilco El Alto'
$ws.Range('B331').Value = 'Autlán De Navarro'
$ws.Range('B338').Value = 'Concepción De Buenos Aires'
$ws.Range('B342').Value = 'Encarnación De Díaz'
$ws.Range('B348').Value = 'Huejuquilla El Alto'
$ws.Range('B349').Value = 'Ixtlahuacán De Los Membrillos'
$ws.Range('B350').Value = 'Ixtlahuacán Del Río'
$ws.Range('B354').Value = 'Jilotlán De Los Dolores'
$ws.Range('B359').Value = 'La Manzanilla De La Paz'
$ws.Range('B360').Value = 'Lagos De Moreno'
$ws.Range('B365').Value = 'Ojuelos De Jalisco'
$ws.Range('B369').Value = 'San Cristóbal De La Barranca'
$ws.Range('B370').Value = 'San Diego De Alejandría'
$ws.Range('B371').Value = 'San Juan De Los Lagos'
$ws.Range('B374').Value = 'San Miguel El Alto'
$ws.Range('B375').Value = 'San Sebastián Del Oeste'
$ws.Range('B381').Value = 'Teocuitatlán De Corona'
$ws.Range('B382').Value = 'Tepatitlán De Morelos'
$ws.Range('B383').Value = 'Tizapán El Alto'
$ws.Range('B384').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B392').Value = 'Unión De San Antonio'
$ws.Range('B393').Value = 'Unión De Tula'
$ws.Range('B394').Value = 'Valle De Guadalupe'
$ws.Range('B395').Value = 'Valle De Juárez'
$ws.Range('B398').Value = 'Zacoalco De Torres'
$ws.Range('B400').Value = 'Zapotitlán De Vadillo'
$ws.Range('B401').Value = 'Zapotlán Del Rey'
$ws.Range('B402').Value = 'Zapotlán El Grande'
$ws.Range('B493').Value = 'Ixtlán Del Río'
$ws.Range('B498').Value = 'Santa María Del Oro'
$ws.Range('B511').Value = 'Chalcatongo De Hidalgo'
$ws.Range('B514').Value = 'Coicoyán De Las Flores'
$ws.Range('B515').Value = 'Constancia Del Rosario'
$ws.Range('B517').Value = 'El Barrio De La Soledad'
$ws.Range('B518').Value = 'Guadalupe De Ramírez'
$ws.Range('B519').Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Range('B520').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B522').Value = 'Ixtlán De Juárez'
$ws.Range('B523').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B527').Value = 'Mariscala De Juárez'
$ws.Range('B528').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B529').Value = 'Oaxaca De Juárez'
$ws.Range('B530').Value = 'Ocotlán De Morelos'
$ws.Range('B531').Value = 'Putla Villa De Guerrero'
$ws.Range('B532').Value = 'San Antonino El Alto'
$ws.Range('B565').Value = 'Santa Cruz Tacache De Mina'
$ws.Range('B568').Value = 'Santa Inés Del Monte'
$ws.Range('B576').Value = 'Santiago Del Río'
$ws.Range('B593').Value = 'Tataltepec De Valdés'
$ws.Range('B594').Value = 'Tlacolula De Matamoros'
$ws.Range('B596').Value = 'Villa Sola De Vega'
$ws.Range('B597').Value = 'Zimatlán De Álvarez'
$ws.Range('B611').Value = 'Cuayuca De Andrade'
$ws.Range('B617').Value = 'Huehuetlán El Chico'
$ws.Range('B618').Value = 'Izúcar De Matamoros'
$ws.Range('B623').Value = 'Palmar De Bravo'
$ws.Range('B626').Value = 'San Salvador El Seco'
$ws.Range('B630').Value = 'Tetela De Ocampo'
$ws.Range('B634').Value = 'Tlacotepec De Benito Juárez'
$ws.Range('B641').Value = 'Amealco De Bonfil'
$ws.Range('B643').Value = 'Cadereyta De Montes'
$ws.Range('B647').Value = 'Jalpan De Serra'
$ws.Range('B650').Value = 'San Juan Del Río'
$ws.Range('B658').Value = 'Ciudad Del Maíz'
$ws.Range('B664').Value = 'Mexquitic De Carmona'
$ws.Range('B668').Value = 'San Ciro De Acosta'
$ws.Range('B671').Value = 'Santa María Del Río'
$ws.Range('B674').Value = 'Villa De Guadalupe'
$ws.Range('B675').Value = 'Villa De La Paz'
$ws.Range('B676').Value = 'Villa De Ramos'
$ws.Range('B677').Value = 'Villa De Reyes'
$ws.Range('B708').Value = 'Nacozari De García'
$ws.Range('B745').Value = 'Tetla De La Solidaridad'
$ws.Range('B752').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B755').Value = 'Castillo De Teayo'
$ws.Range('B764').Value = 'Cosamaloapan De Carpio'
$ws.Range('B768').Value = 'Hueyapan De Ocampo'
$ws.Range('B769').Value = 'Ixhuatlán Del Café'
$ws.Range('B770').Value = 'Ixhuatlán Del Sureste'
$ws.Range('B776').Value = 'Lerdo De Tejada'
$ws.Range('B777').Value = 'Martínez De La Torre'
$ws.Range('B786').Value = 'Paso Del Macho'
$ws.Range('B789').Value = 'Poza Rica De Hidalgo'
$ws.Range('B791').Value = 'Sayula De Alemán'
$ws.Range('B793').Value = 'Soledad De Doblado'
$ws.Range('B799').Value = 'Vega De Alatorre'
$ws.Range('B805').Value = 'Zontecomatlán De López Y Fuentes'
$ws.Range('B822').Value = 'Jiménez Del Teul'
$ws.Range('B826').Value = 'Mezquital Del Oro'
$ws.Range('B830').Value = 'Nochistlán De Mejía'
$ws.Range('B831').Value = 'Noria De Ángeles'
$ws.Range('B841').Value = 'Villa De Cos'

# --- Fix D404 floating point value to match recalculated representation ---
$ws.Range('D404').Value = 0.0904608788853162

# --- Remove trailing footnote rows (848:852), shrinking used range to A1:D846 ---
$ws.Rows("848:852").Delete()
